# Fetch employee status for mystatus page
# Adds a new "Status" column (X) to Sheet1 with a short list of status
# values used to drive the mystatus page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New header cell - reuse the existing header formatting (style index 9)
# from the rest of row 1 (e.g. N1) by copying it across instead of
# re-building the style from scratch.
$ws.Range("X1").Value = "Status"
$ws.Range("N1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New status values for rows 2-7 (plain, unstyled cells like the rest of
# the data rows).
$ws.Range("X2").Value = "Completed/Submitted"
$ws.Range("X3").Value = "Hold/Suspended"
$ws.Range("X4").Value = "ContactManager"
$ws.Range("X5").Value = "Raised a Tickect"
$ws.Range("X6").Value = "Placed for Review"
$ws.Range("X7").Value = "Flipped"

# Match the author's scroll position / selection when they made the edit.
$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X11").Select()
